$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 1119
$ws.Range("F7").Value = 594
$ws.Range("F8").Value = 1522
$ws.Range("F9").Value = 146
$ws.Range("F10").Value = 146
$ws.Range("F11").Value = 1436
$ws.Range("F12").Value = 3062
$ws.Range("F13").Value = 589
$ws.Range("F14").Value = 1737
$ws.Range("F16").Value = 836
$ws.Range("F22").Value = 1190
$ws.Range("F25").Value = 77
$ws.Range("F26").Value = 4703
$ws.Range("F29").Value = 1622
$ws.Range("F30").Value = 57
$ws.Range("F31").Value = 106

$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 47
$ws.Range("F7").Value = 61

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 34

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 34
$ws.Range("F9").Value = 47
$ws.Range("F10").Value = 61
$ws.Range("F12").Value = 1119
$ws.Range("F15").Value = 594
$ws.Range("F16").Value = 1522
$ws.Range("F17").Value = 146
$ws.Range("F18").Value = 146
$ws.Range("F20").Value = 1436
$ws.Range("F21").Value = 3062
$ws.Range("F22").Value = 589
$ws.Range("F23").Value = 1737
$ws.Range("F25").Value = 836
$ws.Range("F33").Value = 1190
$ws.Range("F36").Value = 77
$ws.Range("F37").Value = 4703
$ws.Range("F40").Value = 1622
$ws.Range("F43").Value = 57
$ws.Range("F44").Value = 106

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("I31").Value = "//i2.hdslb.com/bfs/openplatform/202404/QLOhW4Nr1714384036670.png"

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("I44").Value = "//i2.hdslb.com/bfs/openplatform/202404/QLOhW4Nr1714384036670.png"
